# "Generate Report for Archive"
#
# The handoff/localization status for the two tracked files flips from
# "Ready for handoff" to "In Translation" everywhere it is shown:
#   - Overview sheet: the per-language roll-up columns (zh-cn = E, de-de = F)
#   - zh-cn / de-de detail sheets: the "Status" column (C)
#
# The Status column on the detail sheets (and the mirrored columns on the
# Overview sheet) is sized to fit its longest value; "In Translation" is
# shorter than "Ready for handoff", so those columns narrow to match.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: zh-cn (E) / de-de (F) status roll-up columns ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Range("E3").Value = "In Translation"
$wsOverview.Range("F3").Value = "In Translation"
$wsOverview.Range("E1:F1").ColumnWidth = 12.55

# --- zh-cn detail sheet: Status column (C) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Range("C3").Value = "In Translation"
$wsZhCn.Range("C1").ColumnWidth = 12.55

# --- de-de detail sheet: Status column (C) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Range("C3").Value = "In Translation"
$wsDeDe.Range("C1").ColumnWidth = 12.55
